$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V1").ColumnWidth = 28.6667

$ws.Range("V890").Value = "PosControl_ColiGuard"
$ws.Range("V891").Value = "PosControl_Protexin"
$ws.Range("V892").Value = "PosControl_Protexin"
$ws.Range("V893").Value = "PosControl_ColiGuard"
$ws.Range("V894").Value = "PosControl_Protexin"
$ws.Range("V895").Value = "PosControl_ColiGuard"
$ws.Range("V896").Value = "PosControl_ColiGuard"
$ws.Range("V897").Value = "PosControl_Protexin"
$ws.Range("V898").Value = "PosControl_Protexin"
$ws.Range("V899").Value = "PosControl_ColiGuard"
$ws.Range("V900").Value = "PosControl_Protexin"
$ws.Range("V901").Value = "PosControl_ColiGuard"
$ws.Range("V902").Value = "PosControl_Protexin"
$ws.Range("V903").Value = "PosControl_ColiGuard"
$ws.Range("V904").Value = "PosControl_Protexin"
$ws.Range("V905").Value = "PosControl_ColiGuard"

$ws.Range("V896").Select()
